$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data rows (rows 2-7) covering the expanded Sending/Target cluster pairs,
# per "Natmi following Dr Hou advice".
$data = @(
    @("ECs", "Gnai2", "Lhcgr", "FAPs", 3, 1, 203.7816646666667, 611.344994, 0.6667327591988204, 0.6667327591988205, 1, 0.3333333333333333, 0.06641766666666667, 0.199253, 0.01713876626480986, 0.01713876626480986, 13.53470267660911, 121.812324089482, 0.01142697692100034, 0.01142697692100034),
    @("ECs", "Gnai2", "Lhcgr", "sCs", 3, 1, 203.7816646666667, 611.344994, 0.6667327591988204, 0.6667327591988205, 3, 1, 3.808871, 11.426613, 0.9828612337351901, 0.9828612337351901, 776.1780728805913, 6985.602655925322, 0.65530578227782, 0.6553057822778201),
    @("FAPs", "Gnai2", "Lhcgr", "FAPs", 3, 1, 63.14058933333333, 189.421768, 0.2065833519051582, 0.2065833519051582, 1, 0.3333333333333333, 0.06641766666666667, 0.199253, 0.01713876626480986, 0.01713876626480986, 4.193650615478222, 37.742855539304, 0.003540583782503469, 0.00354058378250347),
    @("FAPs", "Gnai2", "Lhcgr", "sCs", 3, 1, 63.14058933333333, 189.421768, 0.2065833519051582, 0.2065833519051582, 3, 1, 3.808871, 11.426613, 0.9828612337351901, 0.9828612337351901, 240.4943596346426, 2164.449236711784, 0.2030427681226547, 0.2030427681226547),
    @("sCs", "Gnai2", "Lhcgr", "FAPs", 3, 1, 38.719942, 116.159826, 0.1266838888960214, 0.1266838888960214, 1, 0.3333333333333333, 0.06641766666666667, 0.199253, 0.01713876626480986, 0.01713876626480986, 2.571688201108667, 23.145193809978, 0.002171205561306052, 0.002171205561306053),
    @("sCs", "Gnai2", "Lhcgr", "sCs", 3, 1, 38.719942, 116.159826, 0.1266838888960214, 0.1266838888960214, 3, 1, 3.808871, 11.426613, 0.9828612337351901, 0.9828612337351901, 147.479264205482, 1327.313377849338, 0.1245126833347153, 0.1245126833347154)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
